$d = $word.ActiveDocument

# The client name "VillageTech" is currently split across two runs ("Village" and
# "Tech"), wrapped in a <w:proofErr w:type="spellStart"/> ... <w:proofErr w:type="spellEnd"/>
# pair (Word's spell-check markers). We need to replace it with a single new run
# containing "Sol Little By Little Enterprise", with the proofErr markers removed,
# while leaving the preceding "| Client: " run completely untouched.
$find = $d.Content
$found = $find.Find.Execute("VillageTech")
if (-not $found) {
    throw "Could not find 'VillageTech' in the document"
}
$vtStart = $find.Start
$vtEnd = $find.End

# --- Left edge -----------------------------------------------------------
# There is nothing between "| Client: " and "Village" other than the
# <w:proofErr spellStart/> marker, and both runs share identical (bold)
# formatting, so any in-place edit right at that seam tends to merge into
# the preceding "| Client: " run. To avoid mutating that run, first insert a
# tiny marker run immediately before "Village" ...
$leftMarker = "QQ"
$leftIns = $d.Range($vtStart, $vtStart)
$leftIns.InsertBefore($leftMarker)

# ... then give just that marker different formatting (non-bold) so Word
# keeps it as its own run instead of silently folding it back into
# "| Client: ". This leaves "| Client: " byte-for-byte identical to before.
$leftMarkerRange = $d.Range($vtStart, $vtStart + $leftMarker.Length)
$leftMarkerRange.Font.Bold = 0

$vtStart = $vtStart + $leftMarker.Length
$vtEnd = $vtEnd + $leftMarker.Length

# --- Right edge ------------------------------------------------------------
# Symmetrically, there is nothing between "Tech" and the paragraph mark other
# than the <w:proofErr spellEnd/> marker, so give it a following run to
# "disappear into" as well.
$rightMarker = "ZZ"
$rightIns = $d.Range($vtEnd, $vtEnd)
$rightIns.InsertAfter($rightMarker)

# --- Replace -----------------------------------------------------------
# Now replace everything from the left marker through the right marker (i.e.
# crossing both former proofErr boundaries) with the new client name in one
# shot. This merges/removes the old "Village"/"Tech" runs (and drops the now
# orphaned proofErr markers with them), producing a single clean new run.
$whole = $d.Range($vtStart - $leftMarker.Length, $vtEnd + $rightMarker.Length)
$whole.Text = "Sol Little By Little Enterprise"
$whole.Font.Bold = 1
